$wb = $excel.ActiveWorkbook

# Existing "ValidLogin" sheet
$wsValid = $wb.Worksheets.Item("ValidLogin")

# Add new worksheet after ValidLogin for the invalid-login scenario
$wsInvalid = $wb.Worksheets.Add($null, $wsValid)
$wsInvalid.Name = "InvalidLogin"

# Populate the InvalidLogin sheet with header + bad credentials
$wsInvalid.Range("A1").Value = "Username"
$wsInvalid.Range("B1").Value = "Password"
$wsInvalid.Range("A2").Value = "abcd"
$wsInvalid.Range("B2").Value = "xyz"

# Selection / active cell bookkeeping to match the target state
$wsValid.Range("A1:B2").Select()
$wsInvalid.Range("B3").Select()

# Make InvalidLogin the active (visible/selected) sheet/tab
$wsInvalid.Activate()
